# Update conditions: the "audioFalse" column becomes "currentPhase",
# and the per-row audio file names in that column are replaced by the
# literal phase identifier "train2P2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header C1: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Data rows C2/C3: file names -> train2P2
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
